$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.369.96'
$ws.Range('E2').Value = '  -3.68%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.848.05'
$ws.Range('E3').Value = '  -5.52%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.65%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.17'
$ws.Range('E5').Value = '  +0.03%  '

# Row 6
$ws.Range('E6').Value = '  -0.60%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4460'
$ws.Range('E7').Value = '  -6.32%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3826'
$ws.Range('E8').Value = '  -5.62%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.97'
$ws.Range('E9').Value = '  -8.48%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07816'
$ws.Range('E10').Value = '  -7.69%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.014'
$ws.Range('E11').Value = '  -4.08%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.39'
$ws.Range('E12').Value = '  -3.57%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.853.11'
$ws.Range('E13').Value = '  -5.26%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.834'
$ws.Range('E14').Value = '  -5.00%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.086'
$ws.Range('E15').Value = '  -6.80%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.005'
$ws.Range('E16').Value = '  -0.55%  '

# Row 17
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '85.14'
$ws.Range('E17').Value = '  -4.89%  '

# Row 18
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001021'
$ws.Range('E18').Value = '  -4.24%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06497'
$ws.Range('E19').Value = '  -1.53%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.92'
$ws.Range('E20').Value = '  -8.95%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.73%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.461'
$ws.Range('E22').Value = '  -5.99%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.380.56'
$ws.Range('E23').Value = '  -3.66%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.75'
$ws.Range('E24').Value = '  -6.64%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.260'
$ws.Range('E25').Value = '  -1.10%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.068.24'
$ws.Range('E26').Value = '  -5.53%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.43'
$ws.Range('E27').Value = '  -2.17%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.30'
$ws.Range('E28').Value = '  -4.45%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.050'
$ws.Range('E29').Value = '  -4.84%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.484'
$ws.Range('E30').Value = '  -7.20%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.86'
$ws.Range('E31').Value = '  -2.99%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09311'
$ws.Range('E32').Value = '  -2.85%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.468'
$ws.Range('E33').Value = '  +1.97%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9254'
$ws.Range('E34').Value = '  -5.29%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.628'
$ws.Range('E35').Value = '  -0.84%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.205'
$ws.Range('E36').Value = '  -6.70%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02214'
$ws.Range('E37').Value = '  -4.77%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05937'
$ws.Range('E38').Value = '  -4.28%  '

# Row 39
$ws.Range('E39').Value = '  -3.59%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.317'
$ws.Range('E40').Value = '  -5.70%  '

# Row 41
$ws.Range('E41').Value = '  -0.66%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5904'
$ws.Range('E42').Value = '  -4.71%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1846'
$ws.Range('E43').Value = '  -3.64%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.25'
$ws.Range('E44').Value = '  -7.81%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.259'
$ws.Range('E45').Value = '  -6.12%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5653'
$ws.Range('E46').Value = '  -5.05%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.17'
$ws.Range('E47').Value = '  -6.18%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.357'
$ws.Range('E48').Value = '  -1.13%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.913'
$ws.Range('E49').Value = '  -6.87%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06862'
$ws.Range('E50').Value = '  +0.88%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '107.93'
$ws.Range('E51').Value = '  -2.26%  '
